$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 1809.8334
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1809.8334
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 5429.5002
$ws.Range("N97").Value = -6421.5002

$ws.Range("H111").Value = 6751.3335
$ws.Range("I111").Value = 9250
$ws.Range("J111").Value = 3628
$ws.Range("K111").Value = 27750
$ws.Range("L111").Value = 10884
$ws.Range("M111").Value = -24683
$ws.Range("N111").Value = -17018

$ws.Range("H112").Value = 1305.9642
$ws.Range("I112").Value = 473.33334
$ws.Range("J112").Value = 1405.88
$ws.Range("K112").Value = 1420.00002
$ws.Range("L112").Value = 4217.64
$ws.Range("M112").Value = -312.0000199999999
$ws.Range("N112").Value = -6433.64

$ws.Range("H118").Value = 1242.5
$ws.Range("I118").Value = 188
$ws.Range("J118").Value = 3000
$ws.Range("K118").Value = 564
$ws.Range("L118").Value = 9000
$ws.Range("M118").Value = 1093
$ws.Range("N118").Value = -12314

$ws.Range("H132").Value = 1850.1875
$ws.Range("I132").Value = 1919.5172
$ws.Range("J132").Value = 1180
$ws.Range("K132").Value = 5758.5516
$ws.Range("L132").Value = 3540
$ws.Range("M132").Value = -3228.5516
$ws.Range("N132").Value = -8600

$ws.Range("H135").Value = 31915520
$ws.Range("I135").Value = 12500658
$ws.Range("J135").Value = 142857580
$ws.Range("K135").Value = 112505922
$ws.Range("L135").Value = 1285718220
$ws.Range("M135").Value = -112503387
$ws.Range("N135").Value = -1285723290

$ws.Range("H137").Value = 1506.5146
$ws.Range("I137").Value = 1106.1632
$ws.Range("J137").Value = 2539
$ws.Range("K137").Value = 3318.4896
$ws.Range("L137").Value = 7617
$ws.Range("M137").Value = -768.4895999999999
$ws.Range("N137").Value = -12717

$ws.Range("H138").Value = 2669.55
$ws.Range("I138").Value = 1203.826
$ws.Range("J138").Value = 4652.5884
$ws.Range("K138").Value = 3611.478
$ws.Range("L138").Value = 13957.7652
$ws.Range("M138").Value = 1528.522
$ws.Range("N138").Value = -24237.7652

$ws.Range("H141").Value = 2323.758
$ws.Range("I141").Value = 1211.1957
$ws.Range("J141").Value = 5522.375
$ws.Range("K141").Value = 3633.5871
$ws.Range("L141").Value = 16567.125
$ws.Range("M141").Value = 1546.4129
$ws.Range("N141").Value = -26927.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23105.281
$ws.Range("I32").Value = 25221.021
$ws.Range("J32").Value = 13161.3
$ws.Range("K32").Value = 25221.021
$ws.Range("L32").Value = 13161.3
$ws.Range("M32").Value = -24934.021
$ws.Range("N32").Value = -13735.3

$ws.Range("H61").Value = 7219.365
$ws.Range("I61").Value = 3888.0952
$ws.Range("J61").Value = 21210.7
$ws.Range("K61").Value = 3888.0952
$ws.Range("L61").Value = 21210.7
$ws.Range("M61").Value = -3676.0952
$ws.Range("N61").Value = -21634.7

$ws.Range("H74").Value = 5187.484
$ws.Range("I74").Value = 2204.0417
$ws.Range("J74").Value = 15416.429
$ws.Range("K74").Value = 2204.0417
$ws.Range("L74").Value = 15416.429
$ws.Range("M74").Value = -1330.0417
$ws.Range("N74").Value = -17164.429

$ws.Range("H77").Value = 5187.484
$ws.Range("I77").Value = 2204.0417
$ws.Range("J77").Value = 15416.429
$ws.Range("K77").Value = 11020.2085
$ws.Range("L77").Value = 77082.145
$ws.Range("M77").Value = -6652.208500000001
$ws.Range("N77").Value = -85818.145

$ws.Range("H97").Value = 1077.619
$ws.Range("I97").Value = 761.3333
$ws.Range("J97").Value = 1868.3334
$ws.Range("K97").Value = 761.3333
$ws.Range("L97").Value = 1868.3334
$ws.Range("M97").Value = -265.3333
$ws.Range("N97").Value = -2860.3334

$ws.Range("H110").Value = 1435.125
$ws.Range("I110").Value = 1359.1111
$ws.Range("J110").Value = 1532.8572
$ws.Range("K110").Value = 1359.1111
$ws.Range("L110").Value = 1532.8572
$ws.Range("M110").Value = 685.8888999999999
$ws.Range("N110").Value = -5622.8572

$ws.Range("H132").Value = 3643.4285
$ws.Range("I132").Value = 3153
$ws.Range("J132").Value = 4869.5
$ws.Range("K132").Value = 9459
$ws.Range("L132").Value = 14608.5
$ws.Range("M132").Value = -6929
$ws.Range("N132").Value = -19668.5

$ws.Range("H135").Value = 27492.111
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 27492.111
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 27492.111
$ws.Range("N135").Value = -37632.111

$ws.Range("H136").Value = 7219.365
$ws.Range("I136").Value = 3888.0952
$ws.Range("J136").Value = 21210.7
$ws.Range("K136").Value = 11664.2856
$ws.Range("L136").Value = 63632.10000000001
$ws.Range("M136").Value = -9114.285600000001
$ws.Range("N136").Value = -68732.10000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1865.0588
$ws.Range("I86").Value = 1733.2667
$ws.Range("J86").Value = 2853.5
$ws.Range("K86").Value = 1733.2667
$ws.Range("L86").Value = 2853.5
$ws.Range("M86").Value = -610.2666999999999
$ws.Range("N86").Value = -5099.5

$ws.Range("H89").Value = 1865.0588
$ws.Range("I89").Value = 1733.2667
$ws.Range("J89").Value = 2853.5
$ws.Range("K89").Value = 8666.333499999999
$ws.Range("L89").Value = 14267.5
$ws.Range("M89").Value = -3050.333499999999
$ws.Range("N89").Value = -25499.5

$ws.Range("H134").Value = 22661.959
$ws.Range("I134").Value = 1573.1111
$ws.Range("J134").Value = 85928.5
$ws.Range("K134").Value = 4719.3333
$ws.Range("L134").Value = 257785.5
$ws.Range("M134").Value = -2184.3333
$ws.Range("N134").Value = -262855.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 5041.25
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 5041.25
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 5041.25
$ws.Range("N13").Value = -5319.25

$ws.Range("H31").Value = 1959.7192
$ws.Range("I31").Value = 1343.9149
$ws.Range("J31").Value = 4854
$ws.Range("K31").Value = 1343.9149
$ws.Range("L31").Value = 4854
$ws.Range("M31").Value = -1048.9149
$ws.Range("N31").Value = -5444

$ws.Range("H34").Value = 1959.7192
$ws.Range("I34").Value = 1343.9149
$ws.Range("J34").Value = 4854
$ws.Range("K34").Value = 1343.9149
$ws.Range("L34").Value = 4854
$ws.Range("M34").Value = -1141.9149
$ws.Range("N34").Value = -5258

$ws.Range("H58").Value = 1110179.5
$ws.Range("I58").Value = 1516209.4
$ws.Range("J58").Value = 2825.4546
$ws.Range("K58").Value = 1516209.4
$ws.Range("L58").Value = 2825.4546
$ws.Range("M58").Value = -1516006.4
$ws.Range("N58").Value = -3231.4546

$ws.Range("H80").Value = 15000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 15000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 15000
$ws.Range("N80").Value = -17246

$ws.Range("H83").Value = 15000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 15000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 45000
$ws.Range("N83").Value = -56232

$ws.Range("H87").Value = 40000
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 40000
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 40000
$ws.Range("N87").Value = -42372

$ws.Range("H90").Value = 40000
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 40000
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 120000
$ws.Range("N90").Value = -131856

$ws.Range("H134").Value = 2909.7234
$ws.Range("I134").Value = 2065.75
$ws.Range("J134").Value = 3534.889
$ws.Range("K134").Value = 6197.25
$ws.Range("L134").Value = 10604.667
$ws.Range("M134").Value = -3662.25
$ws.Range("N134").Value = -15674.667

$ws.Range("H136").Value = 1110179.5
$ws.Range("I136").Value = 1516209.4
$ws.Range("J136").Value = 2825.4546
$ws.Range("K136").Value = 4548628.199999999
$ws.Range("L136").Value = 8476.363799999999
$ws.Range("M136").Value = -4546078.199999999
$ws.Range("N136").Value = -13576.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 20000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 20000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 20000
$ws.Range("N15").Value = -20576

$ws.Range("H81").Value = 20000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 20000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 20000
$ws.Range("N81").Value = -21996

$ws.Range("H84").Value = 20000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 20000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 60000
$ws.Range("N84").Value = -69984

$ws.Range("H132").Value = 7594.0586
$ws.Range("I132").Value = 5120
$ws.Range("J132").Value = 26149.5
$ws.Range("K132").Value = 15360
$ws.Range("L132").Value = 78448.5
$ws.Range("M132").Value = -12830
$ws.Range("N132").Value = -83508.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1023.06665
$ws.Range("I22").Value = 1005.3333
$ws.Range("J22").Value = 1034.8889
$ws.Range("K22").Value = 1005.3333
$ws.Range("L22").Value = 1034.8889
$ws.Range("M22").Value = -710.3333
$ws.Range("N22").Value = -1624.8889

$ws.Range("H27").Value = 1023.06665
$ws.Range("I27").Value = 1005.3333
$ws.Range("J27").Value = 1034.8889
$ws.Range("K27").Value = 1005.3333
$ws.Range("L27").Value = 1034.8889
$ws.Range("M27").Value = -898.3333
$ws.Range("N27").Value = -1248.8889

$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = ""

$ws.Range("H55").Value = 255.78947
$ws.Range("I55").Value = 226
$ws.Range("J55").Value = 296.75
$ws.Range("K55").Value = 226
$ws.Range("L55").Value = 296.75
$ws.Range("M55").Value = -53
$ws.Range("N55").Value = -642.75

$ws.Range("H132").Value = 2199.2222
$ws.Range("I132").Value = 1875.56
$ws.Range("J132").Value = 6245
$ws.Range("K132").Value = 5626.68
$ws.Range("L132").Value = 18735
$ws.Range("M132").Value = -3096.68
$ws.Range("N132").Value = -23795

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = ""
$ws.Range("N14").Value = ""

$ws.Range("H86").Value = 29333.25
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 29333.25
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 29333.25
$ws.Range("N86").Value = -31579.25

$ws.Range("H89").Value = 29333.25
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 29333.25
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 146666.25
$ws.Range("N89").Value = -157898.25

$ws.Range("H93").Value = 68950
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 68950
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 68950
$ws.Range("M93").Value = ""
$ws.Range("N93").Value = -73942

$ws.Range("H132").Value = 807.1622
$ws.Range("I132").Value = 251.34427
$ws.Range("J132").Value = 3415.2307
$ws.Range("K132").Value = 754.0328099999999
$ws.Range("L132").Value = 10245.6921
$ws.Range("M132").Value = 1775.96719
$ws.Range("N132").Value = -15305.6921
